# Verseny ID mezo nem irhato tobbe kezzel: a gomb uj versenyeket general
# automatikus VID_000xx azonositoval. Ez a script az uj, auto-generalt
# versenysorokat (ures start/end/szervezo mezokkel) irja a munkalapra.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 3
$newIds = @("VID_00002", "VID_00003", "VID_00004", "VID_00005", "VID_00006")

for ($i = 0; $i -lt $newIds.Count; $i++) {
    $row = $startRow + $i

    $ws.Cells.Item($row, 1).Value = $newIds[$i]
    $ws.Cells.Item($row, 2).Value = ""
    $ws.Cells.Item($row, 3).Value = ""
    $ws.Cells.Item($row, 4).Value = ""
}
